$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Test 4" block (rows 22-25) used to describe deleting a book ("Vymazanie
# knihy"). Rewrite it to describe the "Add book without a PDF" bug instead,
# keeping the same layout/merged header cell and only touching the text.
$ws.Range("B22").Value = "Test 4: Pridanie knihy"
$ws.Range("C24").Value = "Aplikácia neinformuje používateľa že v knihe chýba PDF."
$ws.Range("C25").Value = "1. Užívateľ klikne na tlačidlo profilu 2. Užívateľ klikne na tlačidlo Add Book                                                     3. Na screene prídávacieho formuláru vyplní všetky polia ale nevyberie PDF.                                     4. Klikne úspešne na tlačidlo Submit 5. Kniha si pridá bez PDF"

# Leave the cursor on the cell that was last edited.
$ws.Range("C25").Select() | Out-Null
